$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "constraints" sentence
# (it starts with a tab, and currently begins with a lastRenderedPageBreak
# immediately followed by the tab and the constraints text).
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "The constraints of the problem are that there are 20 socks",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'constraints' paragraph"
}
$searchRange.Expand(4) | Out-Null   # wdParagraph -> expand to the whole paragraph
$targetRange = $searchRange

# Remove the existing (hidden) _GoBack bookmark; it currently sits right
# after the "My insights..." sentence and needs to move to the end of the
# rewritten "constraints" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$newParagraphXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r>
    <w:tab/>
    <w:t xml:space='preserve'>The constraints of the problem are that there are 20 socks in the drawer of various colors in different amounts. The socks can only be selected in the dark, so </w:t>
  </w:r>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t xml:space='preserve'>determining their color can only be done after selection. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space='preserve'>The sub-goal of this problem is to select the socks without seeing them. </w:t>
  </w:r>
  <w:bookmarkStart w:id='0' w:name='_GoBack'/>
  <w:bookmarkEnd w:id='0'/>
</w:p>
"@

$targetRange.InsertXML($newParagraphXml)
